# Logged Week 17 data for the Raiders' Target Depth Data workbook.
# Updates row 3 ("R" - Road) figures on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 426
$wsOff.Range("C3").Value = 305
$wsOff.Range("D3").Value = 106
$wsOff.Range("E3").Value = 53
$wsOff.Range("F3").Value = 7
$wsOff.Range("G3").Value = 4

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 457
$wsDef.Range("C3").Value = 330
$wsDef.Range("D3").Value = 91
$wsDef.Range("E3").Value = 44
